# Commit: "added p3 pages for produce, added banners for dairy and candy,
# and adjusted banners for mobile view"
#
# The underlying spreadsheet is the team's work-log / contribution tracker
# for the grocery_website project. Tony Yang (row 6) added a new
# contribution note describing the produce-aisle / banner work referenced
# by the commit message, and the active cell moved on to the next empty
# row (C8) ready for the next entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tony Yang's new contribution summary (row 6, column C) - mirrors the
# style (center-aligned) used by the sibling description cells C2/C4.
$ws.Range("C6").Value = "Created template for P2 pages. Made banners for P2 pages. Created Produce aisle (P2), product descriptions for Produce aisle (P3), P5 and P6. Worked on CSS for the pages created."
$ws.Range("C6").HorizontalAlignment = -4108  # xlCenter, matches style used elsewhere in column C

# Move the active selection to C8, matching the saved cursor position.
[void]$ws.Range("C8").Select()
